# Generate Report for Handback
#
# The localization files have come back from handback, so:
#  - Status moves from "Ready for handoff" to "Handed back: in sync with en-US"
#  - The "Latest Target File" (F) and "Latest Handback File" (G) columns get
#    populated with hyperlinks mirroring the existing handoff (A/D) links
#  - The "Latest Handback DateTime" (H) gets the real handback timestamp

$wb = $excel.ActiveWorkbook

$newStatus = "Handed back: in sync with en-US"

# Hyperlink font look (matches the workbook's existing custom "HyperLink"
# style: underlined Calibri in #6495ED) - RGB() order is R + G*256 + B*65536.
$hlColor = 15570276

function Add-HandbackLink($ws, $cellAddr, $text, $url) {
    $ws.Range($cellAddr).Value = $text
    $ws.Hyperlinks.Add($ws.Range($cellAddr), $url, [System.Type]::Missing, [System.Type]::Missing, $text) | Out-Null
    $ws.Range($cellAddr).Style = "HyperLink"
    $ws.Range($cellAddr).Font.Underline = 2
    $ws.Range($cellAddr).Font.Color = $hlColor
}

# ---------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("zh-cn")

$ws.Range("C2").Value = $newStatus
$ws.Range("C3").Value = $newStatus

Add-HandbackLink $ws "F2" "31c3feae-cefe-47da-83cc-2c2c36c01347.md" "https://github.com/OpenLocalizationTest/oltest/blob/7c0f6bc850b9f09065f4dfd281619017ec69fef5/e2e/31c3feae-cefe-47da-83cc-2c2c36c01347.md"
Add-HandbackLink $ws "G2" "31c3feae-cefe-47da-83cc-2c2c36c01347.4150365357c5a12d392aacf644b13ab31a4ecef1.zh-cn.xlf" "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/e3ebc9c5b86842a1474a6ab4e33b3f3272886797/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/31c3feae-cefe-47da-83cc-2c2c36c01347.4150365357c5a12d392aacf644b13ab31a4ecef1.zh-cn.xlf"
Add-HandbackLink $ws "F3" "5d07fbf2-6d78-445b-a87a-800de4aa9db9.md" "https://github.com/OpenLocalizationTest/oltest/blob/7c0f6bc850b9f09065f4dfd281619017ec69fef5/e2e/5d07fbf2-6d78-445b-a87a-800de4aa9db9.md"
Add-HandbackLink $ws "G3" "5d07fbf2-6d78-445b-a87a-800de4aa9db9.46f19e50b5b6a094aeb18b7c88914153792201ce.zh-cn.xlf" "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/e3ebc9c5b86842a1474a6ab4e33b3f3272886797/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/5d07fbf2-6d78-445b-a87a-800de4aa9db9.46f19e50b5b6a094aeb18b7c88914153792201ce.zh-cn.xlf"

$ws.Range("H2").Value = "2016-03-22 02:55:44"
$ws.Range("H3").Value = "2016-03-22 02:55:44"

# ---------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("de-de")

$ws.Range("C2").Value = $newStatus
$ws.Range("C3").Value = $newStatus

Add-HandbackLink $ws "F2" "31c3feae-cefe-47da-83cc-2c2c36c01347.md" "https://github.com/OpenLocalizationTest/oltest/blob/7c0f6bc850b9f09065f4dfd281619017ec69fef5/e2e/31c3feae-cefe-47da-83cc-2c2c36c01347.md"
Add-HandbackLink $ws "G2" "31c3feae-cefe-47da-83cc-2c2c36c01347.4150365357c5a12d392aacf644b13ab31a4ecef1.de-de.xlf" "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/6ced1996db6ff5570caa7c8b1e4555a5227e731b/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/31c3feae-cefe-47da-83cc-2c2c36c01347.4150365357c5a12d392aacf644b13ab31a4ecef1.de-de.xlf"
Add-HandbackLink $ws "F3" "5d07fbf2-6d78-445b-a87a-800de4aa9db9.md" "https://github.com/OpenLocalizationTest/oltest/blob/7c0f6bc850b9f09065f4dfd281619017ec69fef5/e2e/5d07fbf2-6d78-445b-a87a-800de4aa9db9.md"
Add-HandbackLink $ws "G3" "5d07fbf2-6d78-445b-a87a-800de4aa9db9.46f19e50b5b6a094aeb18b7c88914153792201ce.de-de.xlf" "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/6ced1996db6ff5570caa7c8b1e4555a5227e731b/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/5d07fbf2-6d78-445b-a87a-800de4aa9db9.46f19e50b5b6a094aeb18b7c88914153792201ce.de-de.xlf"

$ws.Range("H2").Value = "2016-03-22 02:55:50"
$ws.Range("H3").Value = "2016-03-22 02:55:50"
